$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 512.1
$ws.Range("J103").Value = 707
$ws.Range("L103").Value = 2121
$ws.Range("N103").Value = -3293
$ws.Range("H116").Value = 5051.143
$ws.Range("J116").Value = 6313.3335
$ws.Range("L116").Value = 6313.3335
$ws.Range("N116").Value = -13197.3335
$ws.Range("H135").Value = 1117.0294
$ws.Range("I135").Value = 493
$ws.Range("K135").Value = 4437
$ws.Range("M135").Value = -1902
$ws.Range("H138").Value = 2837.573
$ws.Range("I138").Value = 1420.9062
$ws.Range("J138").Value = 3545.9062
$ws.Range("K138").Value = 4262.7186
$ws.Range("L138").Value = 10637.7186
$ws.Range("M138").Value = 877.2813999999998
$ws.Range("N138").Value = -20917.7186
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2954
$ws.Range("I2").Value = 2362.3333
$ws.Range("J2").Value = 4137.3335
$ws.Range("K2").Value = 2362.3333
$ws.Range("L2").Value = 4137.3335
$ws.Range("M2").Value = -2249.3333
$ws.Range("N2").Value = -4363.3335
$ws.Range("H45").Value = 1644.875
$ws.Range("J45").Value = 1999.75
$ws.Range("L45").Value = 1999.75
$ws.Range("N45").Value = -2753.75
$ws.Range("H97").Value = 859.7
$ws.Range("I97").Value = 561.94116
$ws.Range("K97").Value = 561.94116
$ws.Range("M97").Value = -65.94115999999997
$ws.Range("H110").Value = 3866.8462
$ws.Range("I110").Value = 2808.4443
$ws.Range("J110").Value = 6248.25
$ws.Range("K110").Value = 2808.4443
$ws.Range("L110").Value = 6248.25
$ws.Range("M110").Value = -763.4443000000001
$ws.Range("N110").Value = -10338.25
$ws.Range("H116").Value = 2954
$ws.Range("I116").Value = 2362.3333
$ws.Range("J116").Value = 4137.3335
$ws.Range("K116").Value = 2362.3333
$ws.Range("L116").Value = 4137.3335
$ws.Range("M116").Value = -68.33329999999978
$ws.Range("N116").Value = -8725.333500000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2954
$ws.Range("I3").Value = 2362.3333
$ws.Range("J3").Value = 4137.3335
$ws.Range("K3").Value = 2362.3333
$ws.Range("L3").Value = 4137.3335
$ws.Range("M3").Value = -2248.3333
$ws.Range("N3").Value = -4365.3335
$ws.Range("H86").Value = 2281
$ws.Range("I86").Value = 2226.75
$ws.Range("J86").Value = 2498
$ws.Range("K86").Value = 2226.75
$ws.Range("L86").Value = 2498
$ws.Range("M86").Value = -1103.75
$ws.Range("N86").Value = -4744
$ws.Range("H89").Value = 2281
$ws.Range("I89").Value = 2226.75
$ws.Range("J89").Value = 2498
$ws.Range("K89").Value = 11133.75
$ws.Range("L89").Value = 12490
$ws.Range("M89").Value = -5517.75
$ws.Range("N89").Value = -23722
$ws.Range("H105").Value = 2770
$ws.Range("I105").Value = 2780.1428
$ws.Range("K105").Value = 2780.1428
$ws.Range("M105").Value = -1033.1428
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4774.4634
$ws.Range("I105").Value = 2057.85
$ws.Range("K105").Value = 2057.85
$ws.Range("M105").Value = -310.8499999999999
$ws.Range("H132").Value = 1706.6428
$ws.Range("I132").Value = 1380.1364
$ws.Range("J132").Value = 2903.8333
$ws.Range("K132").Value = 4140.4092
$ws.Range("L132").Value = 8711.499899999999
$ws.Range("M132").Value = -1610.4092
$ws.Range("N132").Value = -13771.4999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 837.125
$ws.Range("J11").Value = 233.33333
$ws.Range("L11").Value = 699.99999
$ws.Range("N11").Value = -979.99999
$ws.Range("H23").Value = 1334.3334
$ws.Range("I23").Value = 310
$ws.Range("J23").Value = 1539.2
$ws.Range("K23").Value = 930
$ws.Range("L23").Value = 4617.6
$ws.Range("M23").Value = -695
$ws.Range("N23").Value = -5087.6
$ws.Range("H34").Value = 1206.1177
$ws.Range("J34").Value = 3747.5
$ws.Range("L34").Value = 11242.5
$ws.Range("N34").Value = -11410.5
$ws.Range("H39").Value = 638.55
$ws.Range("J39").Value = 1996.5
$ws.Range("L39").Value = 5989.5
$ws.Range("N39").Value = -6577.5
$ws.Range("H55").Value = 6523.0586
$ws.Range("J55").Value = 7049.7144
$ws.Range("L55").Value = 21149.1432
$ws.Range("N55").Value = -21503.1432
$ws.Range("H114").Value = 677.625
$ws.Range("I114").Value = 669.8889
$ws.Range("J114").Value = 687.5714
$ws.Range("K114").Value = 2009.6667
$ws.Range("L114").Value = 2062.7142
$ws.Range("M114").Value = 1244.3333
$ws.Range("N114").Value = -8570.7142
$ws.Range("H138").Value = 1814.8462
$ws.Range("I138").Value = 1849.4166
$ws.Range("K138").Value = 5548.2498
$ws.Range("M138").Value = -408.2497999999996
$ws.Range("H140").Value = 4977.3394
$ws.Range("I140").Value = 8208.519
$ws.Range("K140").Value = 24625.557
$ws.Range("M140").Value = -19445.557
$ws.Range("H141").Value = 3799.6
$ws.Range("I141").Value = 3799.6
$ws.Range("K141").Value = 11398.8
$ws.Range("M141").Value = -6218.799999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9702.166999999999
$ws.Range("I70").Value = 9642.6
$ws.Range("K70").Value = 9642.6
$ws.Range("M70").Value = -9372.6
$ws.Range("H73").Value = 9702.166999999999
$ws.Range("I73").Value = 9642.6
$ws.Range("K73").Value = 9642.6
$ws.Range("M73").Value = -8706.6
$ws.Range("H80").Value = 4902.2607
$ws.Range("I80").Value = 4326
$ws.Range("K80").Value = 4326
$ws.Range("M80").Value = -3328
$ws.Range("H83").Value = 4902.2607
$ws.Range("I83").Value = 4326
$ws.Range("K83").Value = 21630
$ws.Range("M83").Value = -16638
$ws.Range("H97").Value = 3582.6924
$ws.Range("I97").Value = 1897.625
$ws.Range("J97").Value = 6278.8
$ws.Range("K97").Value = 1897.625
$ws.Range("L97").Value = 6278.8
$ws.Range("M97").Value = -1401.625
$ws.Range("N97").Value = -7270.8
$ws.Range("H113").Value = 1846
$ws.Range("I113").Value = 1903.6666
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1903.6666
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 266.3334
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 4191.893
$ws.Range("J122").Value = 5503
$ws.Range("L122").Value = 16509
$ws.Range("N122").Value = -21409
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 55874
$ws.Range("J38").Value = 58666.332
$ws.Range("L38").Value = 58666.332
$ws.Range("N38").Value = -59486.332
